# Updated cryptos list data refresh (prices + 1h volume %) and a two-row reorder (Fetch.AI / Aptos).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "56.946.29"
$ws.Range("E2").Value = "  +0.86%  "
$ws.Range("D3").Value = "2.393.45"
$ws.Range("E3").Value = "  +1.46%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'504.13"
$ws.Range("E5").Value = "  -1.67%  "
$ws.Range("D6").Value = "'132.98"
$ws.Range("E6").Value = "  +4.38%  "
$ws.Range("D7").Value = "'0.998"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'0.553"
$ws.Range("E8").Value = "  -0.02%  "
$ws.Range("D9").Value = "2.407.04"
$ws.Range("E9").Value = "  +1.47%  "
$ws.Range("E10").Value = "  +1.09%  "
$ws.Range("E11").Value = "  -1.39%  "
$ws.Range("E12").Value = "  +1.70%  "
$ws.Range("D13").Value = "'4.58"
$ws.Range("E13").Value = "  -4.78%  "
$ws.Range("D14").Value = "2.821.58"
$ws.Range("E14").Value = "  +1.54%  "
$ws.Range("D15").Value = "56.866.06"
$ws.Range("E15").Value = "  +0.87%  "
$ws.Range("D16").Value = "'21.86"
$ws.Range("E16").Value = "  +2.13%  "
$ws.Range("E17").Value = "  +2.28%  "
$ws.Range("D18").Value = "2.395.03"
$ws.Range("E18").Value = "  +1.29%  "
$ws.Range("D19").Value = "'10.21"
$ws.Range("E19").Value = "  -0.59%  "
$ws.Range("D20").Value = "'309.79"
$ws.Range("E20").Value = "  -0.10%  "
$ws.Range("E21").Value = "  -0.36%  "
$ws.Range("D22").Value = "'6.34"
$ws.Range("E22").Value = "  +3.33%  "
$ws.Range("D23").Value = "'5.84"
$ws.Range("E23").Value = "  +0.31%  "
$ws.Range("D24").Value = "'0.999"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "'65.10"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "'0.998"
$ws.Range("E26").Value = "  -0.29%  "
$ws.Range("E27").Value = "  -0.42%  "
$ws.Range("D28").Value = "'0.377"
$ws.Range("E28").Value = "  -3.14%  "
$ws.Range("D29").Value = "'7.40"
$ws.Range("E29").Value = "  +2.87%  "
$ws.Range("D30").Value = "'173.39"
$ws.Range("E30").Value = "  -0.86%  "
$ws.Range("D31").Value = "0.0₃0723"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").Value = "'1.67"
$ws.Range("E32").Value = "  -0.21%  "
$ws.Range("B33").Value = "Aptos"
$ws.Range("C33").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D33").Value = "'5.93"
$ws.Range("E33").Value = "  -3.11%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "'1.12"
$ws.Range("E34").Value = "  -0.07%  "
$ws.Range("E35").Value = "  +0.13%  "
$ws.Range("E36").Value = "  +0.12%  "
$ws.Range("D37").Value = "'17.94"
$ws.Range("E37").Value = "  +1.40%  "
$ws.Range("D38").Value = "'1.19"
$ws.Range("E38").Value = "  +0.55%  "
$ws.Range("D39").Value = "'3.83"
$ws.Range("E39").Value = "  +3.02%  "
$ws.Range("D40").Value = "'36.72"
$ws.Range("E40").Value = "  +2.84%  "
$ws.Range("D41").Value = "'0.805"
$ws.Range("E41").Value = "  -0.19%  "
$ws.Range("E42").Value = "  +1.25%  "
$ws.Range("D43").Value = "'132.66"
$ws.Range("E43").Value = "  +9.34%  "
$ws.Range("D44").Value = "'4.99"
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("E45").Value = "  -0.48%  "
$ws.Range("D46").Value = "'0.567"
$ws.Range("E46").Value = "  +0.01%  "
$ws.Range("D47").Value = "'252.01"
$ws.Range("E47").Value = "  -0.40%  "
$ws.Range("D48").Value = "'0.0910"
$ws.Range("E48").Value = "  +0.57%  "
$ws.Range("D49").Value = "'0.0488"
$ws.Range("E49").Value = "  -0.14%  "
$ws.Range("D50").Value = "'17.06"
$ws.Range("E50").Value = "  +2.85%  "
$ws.Range("E51").Value = "  +1.05%  "
